# Update the "Metadata" sheet (sheet1 / Property-Value table) for the
# rx-pay-tier StructureDefinition: bump the version, refresh the
# publication date, set a real Publisher, replace the empty "Contact"
# rows with a "Jurisdiction" entry, and drop the now-redundant duplicate
# row so the data shifts up by one.
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

$ws1.Range("B3").Value  = "6.0.0"
$ws1.Range("B8").Value  = "2022-01-21T20:46:54+00:00"
$ws1.Range("B9").Value  = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Row 11 was a duplicate "Contact / No display for ContactDetail" row;
# removing it shifts the remaining rows (old Description..Context) up by
# one, turning A1:B21 into A1:B20.
$ws1.Rows.Item(11).Delete()

# Update the root Extension element's Short/Definition text on the
# "Elements" sheet to match the new Title/Description.
$ws2 = $wb.Worksheets.Item("Elements")
$ws2.Range("K2").Value = "Rx Pay Tier"
$ws2.Range("L2").Value = "Customer-specific code for the payment tier of the drug Claim"
